$wb = $excel.ActiveWorkbook

# Rename "Cleansed Customer Data" sheet to "datageneration_sample_spreadshe"
$sheet = $wb.Worksheets.Item("Cleansed Customer Data")
$sheet.Name = "datageneration_sample_spreadshe"

# Activate the "Formulas" sheet (moves the active/selected tab)
$formulasSheet = $wb.Worksheets.Item("Formulas")
$formulasSheet.Activate()
